# Update crypto price/volume figures to the latest scraped snapshot.
# (GitHub Actions scheduled refresh - see commit message.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.630.70"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "3.587.15"
$ws.Range("E3").Value = "  +0.43%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'609.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "'147.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.78%  "
$ws.Range("D7").Value = "3.587.49"
$ws.Range("E7").Value = "  +0.60%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.490"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("D11").Value = "'8.00"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").Value = "4.195.27"
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "'30.10"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").Value = "3.591.43"
$ws.Range("D17").Value = "66.697.21"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").Value = "'11.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("D20").Value = "'6.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").Value = "'15.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.99%  "
$ws.Range("D22").Value = "'432.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("E23").Value = "  +2.57%  "
$ws.Range("D24").Value = "'79.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.59%  "
$ws.Range("D25").Value = "3.727.09"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  -0.77%  "
$ws.Range("D28").Value = "'8.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.79%  "
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("D30").Value = "'2.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "3.581.97"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").Value = "'1.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.27%  "
$ws.Range("E35").Value = "  -2.78%  "
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D38").Value = "'1.72"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.36%  "
$ws.Range("D39").Value = "'5.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("D40").Value = "'173.67"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("D41").Value = "'0.0857"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").Value = "'46.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.50%  "
$ws.Range("D46").Value = "'2.57"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.04%  "
$ws.Range("D47").Value = "'0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  -2.70%  "
$ws.Range("D49").Value = "'25.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.54%  "
$ws.Range("D50").Value = "'23.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.76%  "
$ws.Range("E51").Value = "  +0.95%  "
